# Update column G ("K") values for rows 2-17 on Sheet1.
# This regenerates the "K" (strikeouts) column from new source data,
# replacing the previous "Strike#" based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 4
    3  = 5
    4  = 4
    5  = 3
    6  = 1
    7  = 9
    8  = 6
    9  = 5
    10 = 8
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 3
    16 = 3
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
